# Update "paises" workbook: refreshed COVID country stats + updated timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 01:23"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6385371
$ws.Range("C4").Value = 50127
$ws.Range("D4").Value = 3629692
$ws.Range("E4").Value = 2563667
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 954
$ws.Range("H4").Value = 192012

# Row 5: Brasil
$ws.Range("B5").Value = 4091801
$ws.Range("C5").Value = 45651
$ws.Range("D5").Value = 3278243
$ws.Range("E5").Value = 687974
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 855
$ws.Range("H5").Value = 125584

# Row 9: Colombia
$ws.Range("B9").Value = 650062
$ws.Range("C9").Value = 8488
$ws.Range("D9").Value = 498221
$ws.Range("E9").Value = 130953
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 270
$ws.Range("H9").Value = 20888

# Rows 27/28: Ucrania and Canada swap order (Canada overtakes Ucrania) + refreshed data
$ws.Range("A27").Value = "Canada"
$ws.Range("B27").Value = 131123
$ws.Range("C27").Value = 630
$ws.Range("D27").Value = 115926
$ws.Range("E27").Value = 6056
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 9141

$ws.Range("A28").Value = "Ucrania"
$ws.Range("B28").Value = 130951
$ws.Range("C28").Value = 2723
$ws.Range("D28").Value = 59676
$ws.Range("E28").Value = 68514
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 51
$ws.Range("H28").Value = 2761

# Row 34: Egipto
$ws.Range("B34").Value = 99582
$ws.Range("C34").Value = 157
$ws.Range("D34").Value = 76305
$ws.Range("E34").Value = 17782
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 16
$ws.Range("H34").Value = 5495

# Row 36: Panama
$ws.Range("B36").Value = 95596
$ws.Range("C36").Value = 682
$ws.Range("D36").Value = 68742
$ws.Range("E36").Value = 24791
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 17
$ws.Range("H36").Value = 2063

# Row 43: Guatemala
$ws.Range("B43").Value = 77040
$ws.Range("C43").Value = 682
$ws.Range("D43").Value = 65029
$ws.Range("E43").Value = 9186
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 21
$ws.Range("H43").Value = 2825

# Rows 47/48: Polonia and Japon swap order (Japon overtakes Polonia) + refreshed data
$ws.Range("A47").Value = "Japon"
$ws.Range("B47").Value = 70268
$ws.Range("C47").Value = 669
$ws.Range("D47").Value = 60417
$ws.Range("E47").Value = 8521
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 11
$ws.Range("H47").Value = 1330

$ws.Range("A48").Value = "Polonia"
$ws.Range("B48").Value = 69820
$ws.Range("C48").Value = 691
$ws.Range("D48").Value = 49820
$ws.Range("E48").Value = 17900
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 2100

# Row 54: Nigeria
$ws.Range("B54").Value = 54743
$ws.Range("C54").Value = 156
$ws.Range("D54").Value = 42816
$ws.Range("E54").Value = 10876
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 3
$ws.Range("H54").Value = 1051

# Row 55: Barein
$ws.Range("B55").Value = 54095
$ws.Range("C55").Value = 662
$ws.Range("D55").Value = 50323
$ws.Range("E55").Value = 3577
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 195

# Row 72: Chequia
$ws.Range("B72").Value = 27249
$ws.Range("C72").Value = 797
$ws.Range("D72").Value = 19027
$ws.Range("E72").Value = 7793
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 429

# Row 83: Bulgaria
$ws.Range("B83").Value = 16954
$ws.Range("C83").Value = 179
$ws.Range("D83").Value = 12046
$ws.Range("E83").Value = 4243
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 7
$ws.Range("H83").Value = 665

# Rows 91/92: Grecia and Noruega swap order (Noruega overtakes Grecia) + refreshed data
$ws.Range("A91").Value = "Noruega"
$ws.Range("B91").Value = 11231
$ws.Range("C91").Value = 111
$ws.Range("D91").Value = 9348
$ws.Range("E91").Value = 1619
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 264

$ws.Range("A92").Value = "Grecia"
$ws.Range("B92").Value = 11200
$ws.Range("C92").Value = 202
$ws.Range("D92").Value = 3804
$ws.Range("E92").Value = 7117
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 279

# Row 99: Gabon
$ws.Range("B99").Value = 8601
$ws.Range("C99").Value = 63
$ws.Range("D99").Value = 7424
$ws.Range("E99").Value = 1124
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 53

# Row 110: Montenegro
$ws.Range("B110").Value = 5275
$ws.Range("C110").Value = 110
$ws.Range("D110").Value = 4203
$ws.Range("E110").Value = 966
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 106

# Row 114: Republica de Africa Central
$ws.Range("B114").Value = 4729
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 1818
$ws.Range("E114").Value = 2849
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 62

# Row 117: Tunez
$ws.Range("B117").Value = 4542
$ws.Range("C117").Value = 148
$ws.Range("D117").Value = 1699
$ws.Range("E117").Value = 2756
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 3
$ws.Range("H117").Value = 87

# Row 122: Surinam
$ws.Range("B122").Value = 4252
$ws.Range("C122").Value = 37
$ws.Range("D122").Value = 3366
$ws.Range("E122").Value = 811
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 2
$ws.Range("H122").Value = 75

# Row 145: Trinidad yTobago
$ws.Range("B145").Value = 2040
$ws.Range("C145").Value = 56
$ws.Range("D145").Value = 707
$ws.Range("E145").Value = 1302
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 31

# Row 151: Uruguay
$ws.Range("B151").Value = 1653
$ws.Range("C151").Value = 17
$ws.Range("D151").Value = 1446
$ws.Range("E151").Value = 162
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 45

# Rows 214/215: Montserrat and Islas Malvinas swap order (Islas Malvinas overtakes Montserrat) + refreshed data
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
